$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.404.12'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.609.42'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.27'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  +2.00%  '
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.070.19'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '59.319.82'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.56'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.661.46'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '343.99'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.36'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.14'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.39'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.408'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.21'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0740'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.70'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +8.48%  '
$ws.Range('E31').Value = '  -1.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.79'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '149.47'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.99'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '37.09'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.00%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.11'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.837'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.811'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.56'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '275.88'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.74'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.952.90'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.31'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.42'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.90%  '
